$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (102) down into the
# three new rows (103-105) so the new rows inherit the same styles
# (bold/bordered index column, date-formatted match-time column, etc.).
$ws.Range("A102:V102").Copy($ws.Range("A103:V105"))

# --- Row 103: PSV 4 x 0 Zwolle ---
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = "netherlands"
$ws.Range("C103").Value = "eredivisie"
$ws.Range("D103").Value = "2023-2024"
$ws.Range("E103").Value = 45242.51041666666
$ws.Range("F103").Value = "PSV"
$ws.Range("G103").Value = 4
$ws.Range("H103").Value = "Zwolle"
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1.09
$ws.Range("K103").Value = "05/11/2023 14:42"
$ws.Range("L103").Value = 1.07
$ws.Range("M103").Value = "12/11/2023 11:38"
$ws.Range("N103").Value = 12.69
$ws.Range("O103").Value = "05/11/2023 14:42"
$ws.Range("P103").Value = 15.68
$ws.Range("Q103").Value = "12/11/2023 12:14"
$ws.Range("R103").Value = 18.16
$ws.Range("S103").Value = "05/11/2023 14:42"
$ws.Range("T103").Value = 32
$ws.Range("U103").Value = "12/11/2023 12:14"
$ws.Range("V103").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/psv-zwolle/WYBTpYz3/"

# --- Row 104: Almere City 2 x 2 Ajax ---
$ws.Range("A104").Value = 103
$ws.Range("B104").Value = "netherlands"
$ws.Range("C104").Value = "eredivisie"
$ws.Range("D104").Value = "2023-2024"
$ws.Range("E104").Value = 45242.60416666666
$ws.Range("F104").Value = "Almere City"
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = "Ajax"
$ws.Range("I104").Value = 2
$ws.Range("J104").Value = 4.82
$ws.Range("K104").Value = "06/11/2023 07:04"
$ws.Range("L104").Value = 5.52
$ws.Range("M104").Value = "12/11/2023 14:27"
$ws.Range("N104").Value = 4.65
$ws.Range("O104").Value = "06/11/2023 07:04"
$ws.Range("P104").Value = 4.94
$ws.Range("Q104").Value = "12/11/2023 14:23"
$ws.Range("R104").Value = 1.61
$ws.Range("S104").Value = "06/11/2023 07:04"
$ws.Range("T104").Value = 1.54
$ws.Range("U104").Value = "12/11/2023 14:21"
$ws.Range("V104").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/almere-city-ajax/KdT5Rqzr/"

# --- Row 105: FC Volendam 1 x 4 Sparta Rotterdam ---
$ws.Range("A105").Value = 104
$ws.Range("B105").Value = "netherlands"
$ws.Range("C105").Value = "eredivisie"
$ws.Range("D105").Value = "2023-2024"
$ws.Range("E105").Value = 45242.60416666666
$ws.Range("F105").Value = "FC Volendam"
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = "Sparta Rotterdam"
$ws.Range("I105").Value = 4
$ws.Range("J105").Value = 3.3
$ws.Range("K105").Value = "05/11/2023 17:12"
$ws.Range("L105").Value = 3.16
$ws.Range("M105").Value = "12/11/2023 14:25"
$ws.Range("N105").Value = 3.84
$ws.Range("O105").Value = "05/11/2023 17:12"
$ws.Range("P105").Value = 3.83
$ws.Range("Q105").Value = "12/11/2023 14:25"
$ws.Range("R105").Value = 2.1
$ws.Range("S105").Value = "05/11/2023 17:12"
$ws.Range("T105").Value = 2.22
$ws.Range("U105").Value = "12/11/2023 14:25"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/fc-volendam-sparta-rotterdam/vaFyqC5F/"
